$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Typography) edits ---
$ws1.Cells.Item(6,3).Value = 'Asap-Regular.ttf'
$ws1.Cells.Item(7,2).Value = 'Login'
$ws1.Cells.Item(7,3).Value = 'Asap-Regular.ttf'
$ws1.Cells.Item(9,2).Value = 'Logo_M'
$ws1.Cells.Item(9,3).Value = 'Asap-Regular.ttf'
$ws1.Cells.Item(9,6).Value = '?'
$ws1.Cells.Item(10,2).Value = 'AI'
$ws1.Cells.Item(10,3).Value = 'Asap-Regular.ttf'
$ws1.Cells.Item(10,6).Value = '?'
$ws1.Cells.Item(11,2).Value = 'b_DI'
$ws1.Cells.Item(11,3).Value = 'Asap-Regular.ttf'
$ws1.Cells.Item(11,6).Value = '?'
$ws1.Cells.Item(12,2).Value = 'b_Screens'
$ws1.Cells.Item(12,3).Value = 'arial.ttf'
$ws1.Cells.Item(12,6).Value = '?'
$ws1.Cells.Item(4,4).Value = 60
$ws1.Cells.Item(7,4).Value = 20
$ws1.Cells.Item(9,4).Value = 30
$ws1.Cells.Item(9,5).Value = 4
$ws1.Cells.Item(10,4).Value = 25
$ws1.Cells.Item(10,5).Value = 4
$ws1.Cells.Item(11,4).Value = 20
$ws1.Cells.Item(11,5).Value = 4
$ws1.Cells.Item(12,4).Value = 12
$ws1.Cells.Item(12,5).Value = 4

# --- Sheet2 (Translation) edits ---
$ws2.Cells.Item(23,2).Value = 'SingleUseId79'
$ws2.Cells.Item(23,3).Value = 'Typography_01'
$ws2.Cells.Item(23,4).Value = 'Left'
$ws2.Cells.Item(23,5).Value = '<>'
$ws2.Cells.Item(23,6).Value = 'LTR'
$ws2.Cells.Item(24,2).Value = 'SingleUseId80'
$ws2.Cells.Item(24,3).Value = 'Typography_01'
$ws2.Cells.Item(24,4).Value = 'Left'
$ws2.Cells.Item(24,5).Value = 'DD/MM/YYYY'
$ws2.Cells.Item(24,6).Value = 'LTR'
$ws2.Cells.Item(25,2).Value = 'SingleUseId81'
$ws2.Cells.Item(25,3).Value = 'Typography_01'
$ws2.Cells.Item(25,4).Value = 'Left'
$ws2.Cells.Item(25,5).Value = '<value>; '
$ws2.Cells.Item(25,6).Value = 'LTR'
$ws2.Cells.Item(26,2).Value = 'SingleUseId82'
$ws2.Cells.Item(26,3).Value = 'Typography_01'
$ws2.Cells.Item(26,4).Value = 'Left'
$ws2.Cells.Item(26,5).Value = '0'
$ws2.Cells.Item(26,6).Value = 'LTR'
$ws2.Cells.Item(27,2).Value = 'SingleUseId83'
$ws2.Cells.Item(27,3).Value = 'Login'
$ws2.Cells.Item(27,4).Value = 'Left'
$ws2.Cells.Item(27,5).Value = 'Login'
$ws2.Cells.Item(27,6).Value = 'LTR'
$ws2.Cells.Item(28,2).Value = 'SingleUseId84'
$ws2.Cells.Item(28,3).Value = 'Login'
$ws2.Cells.Item(28,4).Value = 'Left'
$ws2.Cells.Item(28,5).Value = 'Password'
$ws2.Cells.Item(28,6).Value = 'LTR'
$ws2.Cells.Item(29,2).Value = 'SingleUseId85'
$ws2.Cells.Item(29,3).Value = 'Default'
$ws2.Cells.Item(29,4).Value = 'Left'
$ws2.Cells.Item(29,5).Value = 'ELHART'
$ws2.Cells.Item(29,6).Value = 'LTR'
$ws2.Cells.Item(30,2).Value = 'SingleUseId86'
$ws2.Cells.Item(30,3).Value = 'Logo_M'
$ws2.Cells.Item(30,4).Value = 'Left'
$ws2.Cells.Item(30,5).Value = 'ELHART'
$ws2.Cells.Item(30,6).Value = 'LTR'
$ws2.Cells.Item(31,2).Value = 'SingleUseId87'
$ws2.Cells.Item(31,3).Value = 'Logo_M'
$ws2.Cells.Item(31,4).Value = 'Left'
$ws2.Cells.Item(31,5).Value = 'ELHART'
$ws2.Cells.Item(31,6).Value = 'LTR'
$ws2.Cells.Item(32,2).Value = 'SingleUseId88'
$ws2.Cells.Item(32,3).Value = 'AI'
$ws2.Cells.Item(32,4).Value = 'Left'
$ws2.Cells.Item(32,5).Value = 'AI_1'
$ws2.Cells.Item(32,6).Value = 'LTR'
$ws2.Cells.Item(33,2).Value = 'SingleUseId89'
$ws2.Cells.Item(33,3).Value = 'AI'
$ws2.Cells.Item(33,4).Value = 'Left'
$ws2.Cells.Item(33,5).Value = 'AI_2'
$ws2.Cells.Item(33,6).Value = 'LTR'
$ws2.Cells.Item(34,2).Value = 'SingleUseId90'
$ws2.Cells.Item(34,3).Value = 'AI'
$ws2.Cells.Item(34,4).Value = 'Left'
$ws2.Cells.Item(34,5).Value = 'AI_3'
$ws2.Cells.Item(34,6).Value = 'LTR'
$ws2.Cells.Item(35,2).Value = 'SingleUseId91'
$ws2.Cells.Item(35,3).Value = 'AI'
$ws2.Cells.Item(35,4).Value = 'Left'
$ws2.Cells.Item(35,5).Value = 'AI_4'
$ws2.Cells.Item(35,6).Value = 'LTR'
$ws2.Cells.Item(36,2).Value = 'SingleUseId92'
$ws2.Cells.Item(36,3).Value = 'b_DI'
$ws2.Cells.Item(36,4).Value = 'Center'
$ws2.Cells.Item(36,5).Value = 'DI_1'
$ws2.Cells.Item(36,6).Value = 'LTR'
$ws2.Cells.Item(37,2).Value = 'SingleUseId93'
$ws2.Cells.Item(37,3).Value = 'b_DI'
$ws2.Cells.Item(37,4).Value = 'Center'
$ws2.Cells.Item(37,5).Value = 'DI_2'
$ws2.Cells.Item(37,6).Value = 'LTR'
$ws2.Cells.Item(38,2).Value = 'SingleUseId94'
$ws2.Cells.Item(38,3).Value = 'b_DI'
$ws2.Cells.Item(38,4).Value = 'Center'
$ws2.Cells.Item(38,5).Value = 'DI_3'
$ws2.Cells.Item(38,6).Value = 'LTR'
$ws2.Cells.Item(39,2).Value = 'SingleUseId95'
$ws2.Cells.Item(39,3).Value = 'b_DI'
$ws2.Cells.Item(39,4).Value = 'Center'
$ws2.Cells.Item(39,5).Value = 'DI_4'
$ws2.Cells.Item(39,6).Value = 'LTR'
$ws2.Cells.Item(40,2).Value = 'SingleUseId96'
$ws2.Cells.Item(40,3).Value = 'b_DI'
$ws2.Cells.Item(40,4).Value = 'Center'
$ws2.Cells.Item(40,5).Value = 'DO_1'
$ws2.Cells.Item(40,6).Value = 'LTR'
$ws2.Cells.Item(41,2).Value = 'SingleUseId97'
$ws2.Cells.Item(41,3).Value = 'b_DI'
$ws2.Cells.Item(41,4).Value = 'Center'
$ws2.Cells.Item(41,5).Value = 'DO_2'
$ws2.Cells.Item(41,6).Value = 'LTR'
$ws2.Cells.Item(42,2).Value = 'SingleUseId98'
$ws2.Cells.Item(42,3).Value = 'b_DI'
$ws2.Cells.Item(42,4).Value = 'Center'
$ws2.Cells.Item(42,5).Value = 'DO_3'
$ws2.Cells.Item(42,6).Value = 'LTR'
$ws2.Cells.Item(43,2).Value = 'SingleUseId99'
$ws2.Cells.Item(43,3).Value = 'b_DI'
$ws2.Cells.Item(43,4).Value = 'Center'
$ws2.Cells.Item(43,5).Value = 'DO_4'
$ws2.Cells.Item(43,6).Value = 'LTR'
$ws2.Cells.Item(44,2).Value = 'SingleUseId100'
$ws2.Cells.Item(44,3).Value = 'b_DI'
$ws2.Cells.Item(44,4).Value = 'Center'
$ws2.Cells.Item(44,5).Value = 'DI_1'
$ws2.Cells.Item(44,6).Value = 'LTR'
$ws2.Cells.Item(45,2).Value = 'SingleUseId101'
$ws2.Cells.Item(45,3).Value = 'b_DI'
$ws2.Cells.Item(45,4).Value = 'Center'
$ws2.Cells.Item(45,5).Value = 'DI_2'
$ws2.Cells.Item(45,6).Value = 'LTR'
$ws2.Cells.Item(46,2).Value = 'SingleUseId102'
$ws2.Cells.Item(46,3).Value = 'b_DI'
$ws2.Cells.Item(46,4).Value = 'Center'
$ws2.Cells.Item(46,5).Value = 'DI_3'
$ws2.Cells.Item(46,6).Value = 'LTR'
$ws2.Cells.Item(47,2).Value = 'SingleUseId103'
$ws2.Cells.Item(47,3).Value = 'b_DI'
$ws2.Cells.Item(47,4).Value = 'Center'
$ws2.Cells.Item(47,5).Value = 'DI_4'
$ws2.Cells.Item(47,6).Value = 'LTR'
$ws2.Cells.Item(48,2).Value = 'SingleUseId104'
$ws2.Cells.Item(48,3).Value = 'b_DI'
$ws2.Cells.Item(48,4).Value = 'Center'
$ws2.Cells.Item(48,5).Value = 'DO_1'
$ws2.Cells.Item(48,6).Value = 'LTR'
$ws2.Cells.Item(49,2).Value = 'SingleUseId105'
$ws2.Cells.Item(49,3).Value = 'b_DI'
$ws2.Cells.Item(49,4).Value = 'Center'
$ws2.Cells.Item(49,5).Value = 'DO_2'
$ws2.Cells.Item(49,6).Value = 'LTR'
$ws2.Cells.Item(50,2).Value = 'SingleUseId106'
$ws2.Cells.Item(50,3).Value = 'b_DI'
$ws2.Cells.Item(50,4).Value = 'Center'
$ws2.Cells.Item(50,5).Value = 'DO_3'
$ws2.Cells.Item(50,6).Value = 'LTR'
$ws2.Cells.Item(51,2).Value = 'SingleUseId107'
$ws2.Cells.Item(51,3).Value = 'b_DI'
$ws2.Cells.Item(51,4).Value = 'Center'
$ws2.Cells.Item(51,5).Value = 'DO_4'
$ws2.Cells.Item(51,6).Value = 'LTR'
$ws2.Cells.Item(52,2).Value = 'SingleUseId108'
$ws2.Cells.Item(52,3).Value = 'b_Screens'
$ws2.Cells.Item(52,4).Value = 'Center'
$ws2.Cells.Item(52,5).Value = 'Login'
$ws2.Cells.Item(52,6).Value = 'LTR'
$ws2.Cells.Item(53,2).Value = 'SingleUseId109'
$ws2.Cells.Item(53,3).Value = 'b_Screens'
$ws2.Cells.Item(53,4).Value = 'Center'
$ws2.Cells.Item(53,5).Value = 'Set Points'
$ws2.Cells.Item(53,6).Value = 'LTR'
$ws2.Cells.Item(54,2).Value = 'SingleUseId110'
$ws2.Cells.Item(54,3).Value = 'b_Screens'
$ws2.Cells.Item(54,4).Value = 'Center'
$ws2.Cells.Item(54,5).Value = 'Graphic'
$ws2.Cells.Item(54,6).Value = 'LTR'
$ws2.Cells.Item(55,2).Value = 'SingleUseId111'
$ws2.Cells.Item(55,3).Value = 'b_Screens'
$ws2.Cells.Item(55,4).Value = 'Center'
$ws2.Cells.Item(55,5).Value = 'Analytics'
$ws2.Cells.Item(55,6).Value = 'LTR'
$ws2.Cells.Item(56,2).Value = 'SingleUseId112'
$ws2.Cells.Item(56,3).Value = 'b_Screens'
$ws2.Cells.Item(56,4).Value = 'Center'
$ws2.Cells.Item(56,5).Value = 'Archive'
$ws2.Cells.Item(56,6).Value = 'LTR'
$ws2.Cells.Item(57,2).Value = 'SingleUseId113'
$ws2.Cells.Item(57,3).Value = 'b_Screens'
$ws2.Cells.Item(57,4).Value = 'Center'
$ws2.Cells.Item(57,5).Value = 'Login'
$ws2.Cells.Item(57,6).Value = 'LTR'
$ws2.Cells.Item(58,2).Value = 'SingleUseId114'
$ws2.Cells.Item(58,3).Value = 'b_Screens'
$ws2.Cells.Item(58,4).Value = 'Center'
$ws2.Cells.Item(58,5).Value = 'Set Points'
$ws2.Cells.Item(58,6).Value = 'LTR'
$ws2.Cells.Item(59,2).Value = 'SingleUseId115'
$ws2.Cells.Item(59,3).Value = 'b_Screens'
$ws2.Cells.Item(59,4).Value = 'Center'
$ws2.Cells.Item(59,5).Value = 'Discret'
$ws2.Cells.Item(59,6).Value = 'LTR'
$ws2.Cells.Item(60,2).Value = 'SingleUseId116'
$ws2.Cells.Item(60,3).Value = 'b_Screens'
$ws2.Cells.Item(60,4).Value = 'Center'
$ws2.Cells.Item(60,5).Value = 'Analytics'
$ws2.Cells.Item(60,6).Value = 'LTR'
$ws2.Cells.Item(61,2).Value = 'SingleUseId117'
$ws2.Cells.Item(61,3).Value = 'b_Screens'
$ws2.Cells.Item(61,4).Value = 'Center'
$ws2.Cells.Item(61,5).Value = 'Archive'
$ws2.Cells.Item(61,6).Value = 'LTR'
$ws2.Cells.Item(62,2).Value = 'SingleUseId118'
$ws2.Cells.Item(62,3).Value = 'Typography_01'
$ws2.Cells.Item(62,4).Value = 'Left'
$ws2.Cells.Item(62,5).Value = 'Analog inputs #1'
$ws2.Cells.Item(62,6).Value = 'LTR'
$ws2.Cells.Item(63,2).Value = 'SingleUseId120'
$ws2.Cells.Item(63,3).Value = 'Typography_01'
$ws2.Cells.Item(63,4).Value = 'Left'
$ws2.Cells.Item(63,5).Value = 'dd/mm/yyyy'
$ws2.Cells.Item(63,6).Value = 'LTR'
